$d = $word.ActiveDocument

# The document's headers/footers each contain one inline picture whose
# docPr/cNvPr "name" attribute needs to be swapped:
#   image1.png -> image2.png   (Pearson logo, found in the footers)
#   image2.jpg -> image1.jpg   (BTEC logo,    found in the headers)
$rename = @{
    "image1.png" = "image2.png"
    "image2.jpg" = "image1.jpg"
}

function Rename-InlinePicture($inlineShape) {
    # InlineShape has no writable Name property in the Word object model,
    # but converting to a floating Shape exposes Name for read/write; then
    # converting back restores the original wp:inline (non-anchored) markup.
    $shape = $inlineShape.ConvertToShape()
    $oldName = $shape.Name
    if ($rename.ContainsKey($oldName)) {
        $shape.Name = $rename[$oldName]
    }
    $shape.ConvertToInlineShape() | Out-Null
}

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $sec = $d.Sections.Item($secIdx)

    for ($i = 1; $i -le 3; $i++) {
        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
                Rename-InlinePicture $h.Range.InlineShapes.Item($j)
            }
        }
    }

    for ($i = 1; $i -le 3; $i++) {
        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
                Rename-InlinePicture $f.Range.InlineShapes.Item($j)
            }
        }
    }
}
